# Re-upload of the peer-grading form export: bump the cached "now()" response
# timestamp forward by one month and leave the selection where the user last
# clicked (N4, just past the last data row) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# N3 holds the ReceivedAtUTC value produced by utcNow() for the second
# response row; shift it from Feb 19 to Mar 19 (same time-of-day/ticks).
$ws.Range("N3").Value = "2026-03-19T14:35:52.5219332Z"

# Match the saved cursor/selection position recorded in the workbook.
$ws.Range("N4").Select()
